# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Atomos_Profits workbook sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 18336.5
$ws.Range("I21").Value = 18336.5
$ws.Range("K21").Value = 18336.5
$ws.Range("M21").Value = -17868.5
$ws.Range("H23").Value = 18336.5
$ws.Range("I23").Value = 18336.5
$ws.Range("K23").Value = 18336.5
$ws.Range("M23").Value = -18102.5
$ws.Range("H39").Value = 722.7
$ws.Range("I39").Value = 71.36364
$ws.Range("J39").Value = 1518.7778
$ws.Range("K39").Value = 214.09092
$ws.Range("L39").Value = 4556.3334
$ws.Range("M39").Value = 81.90907999999999
$ws.Range("N39").Value = -5148.3334
$ws.Range("H116").Value = 2991.9487
$ws.Range("I116").Value = 2394.087
$ws.Range("J116").Value = 3851.375
$ws.Range("K116").Value = 2394.087
$ws.Range("L116").Value = 3851.375
$ws.Range("M116").Value = 1047.913
$ws.Range("N116").Value = -10735.375
$ws.Range("H137").Value = 2780559.8
$ws.Range("I137").Value = 7697030.5
$ws.Range("J137").Value = 1684.9131
$ws.Range("K137").Value = 23091091.5
$ws.Range("L137").Value = 5054.7393
$ws.Range("M137").Value = -23088541.5
$ws.Range("N137").Value = -10154.7393
$ws.Range("H138").Value = 4392.7
$ws.Range("I138").Value = 2482.25
$ws.Range("J138").Value = 7258.375
$ws.Range("K138").Value = 7446.75
$ws.Range("L138").Value = 21775.125
$ws.Range("M138").Value = -2306.75
$ws.Range("N138").Value = -32055.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7813877
$ws.Range("I2").Value = 25000910
$ws.Range("J2").Value = 1589.5
$ws.Range("K2").Value = 25000910
$ws.Range("L2").Value = 1589.5
$ws.Range("M2").Value = -25000797
$ws.Range("N2").Value = -1815.5
$ws.Range("H45").Value = 1720.2424
$ws.Range("I45").Value = 1062.96
$ws.Range("J45").Value = 3774.25
$ws.Range("K45").Value = 1062.96
$ws.Range("L45").Value = 3774.25
$ws.Range("M45").Value = -685.96
$ws.Range("N45").Value = -4528.25
$ws.Range("H74").Value = 989
$ws.Range("I74").Value = 942.4
$ws.Range("K74").Value = 942.4
$ws.Range("M74").Value = -68.39999999999998
$ws.Range("H77").Value = 989
$ws.Range("I77").Value = 942.4
$ws.Range("K77").Value = 4712
$ws.Range("M77").Value = -344
$ws.Range("H116").Value = 7813877
$ws.Range("I116").Value = 25000910
$ws.Range("J116").Value = 1589.5
$ws.Range("K116").Value = 25000910
$ws.Range("L116").Value = 1589.5
$ws.Range("M116").Value = -24998616
$ws.Range("N116").Value = -6177.5
$ws.Range("H132").Value = 15628381
$ws.Range("I132").Value = 22225682
$ws.Range("J132").Value = 3194.6316
$ws.Range("K132").Value = 66677046
$ws.Range("L132").Value = 9583.8948
$ws.Range("M132").Value = -66674516
$ws.Range("N132").Value = -14643.8948

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7813877
$ws.Range("I3").Value = 25000910
$ws.Range("J3").Value = 1589.5
$ws.Range("K3").Value = 25000910
$ws.Range("L3").Value = 1589.5
$ws.Range("M3").Value = -25000796
$ws.Range("N3").Value = -1817.5
$ws.Range("H25").Value = 37458
$ws.Range("I25").Value = 4900
$ws.Range("K25").Value = 4900
$ws.Range("M25").Value = -4665
$ws.Range("H75").Value = 2100
$ws.Range("I75").Value = 2100
$ws.Range("K75").Value = 2100
$ws.Range("M75").Value = -1164
$ws.Range("H78").Value = 2100
$ws.Range("I78").Value = 2100
$ws.Range("K78").Value = 6300
$ws.Range("M78").Value = -1620
$ws.Range("H97").Value = 14868.917
$ws.Range("I97").Value = 5685.4
$ws.Range("K97").Value = 5685.4
$ws.Range("M97").Value = -4694.4
$ws.Range("H137").Value = 31709
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = $null
$ws.Range("N140").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1589678
$ws.Range("I31").Value = 3031815.2
$ws.Range("J31").Value = 3327.1
$ws.Range("K31").Value = 3031815.2
$ws.Range("L31").Value = 3327.1
$ws.Range("M31").Value = -3031520.2
$ws.Range("N31").Value = -3917.1
$ws.Range("H34").Value = 1589678
$ws.Range("I34").Value = 3031815.2
$ws.Range("J34").Value = 3327.1
$ws.Range("K34").Value = 3031815.2
$ws.Range("L34").Value = 3327.1
$ws.Range("M34").Value = -3031613.2
$ws.Range("N34").Value = -3731.1
$ws.Range("H58").Value = 31254376
$ws.Range("I58").Value = 3010
$ws.Range("J58").Value = 83339980
$ws.Range("K58").Value = 3010
$ws.Range("L58").Value = 83339980
$ws.Range("M58").Value = -2807
$ws.Range("N58").Value = -83340386
$ws.Range("H112").Value = 27500
$ws.Range("J112").Value = 27500
$ws.Range("L112").Value = 27500
$ws.Range("N112").Value = -30454
$ws.Range("H122").Value = 4336.4546
$ws.Range("I122").Value = 3800.3333
$ws.Range("J122").Value = 4979.8
$ws.Range("K122").Value = 11400.9999
$ws.Range("L122").Value = 14939.4
$ws.Range("M122").Value = -8950.999899999999
$ws.Range("N122").Value = -19839.4
$ws.Range("H132").Value = 4524.9165
$ws.Range("I132").Value = 2966.6667
$ws.Range("J132").Value = 5459.8667
$ws.Range("K132").Value = 8900.000100000001
$ws.Range("L132").Value = 16379.6001
$ws.Range("M132").Value = -6370.000100000001
$ws.Range("N132").Value = -21439.6001
$ws.Range("H134").Value = 1662.18
$ws.Range("I134").Value = 1415.3256
$ws.Range("J134").Value = 3178.5715
$ws.Range("K134").Value = 4245.976799999999
$ws.Range("L134").Value = 9535.7145
$ws.Range("M134").Value = -1710.976799999999
$ws.Range("N134").Value = -14605.7145
$ws.Range("H136").Value = 31254376
$ws.Range("I136").Value = 3010
$ws.Range("J136").Value = 83339980
$ws.Range("K136").Value = 9030
$ws.Range("L136").Value = 250019940
$ws.Range("M136").Value = -6480
$ws.Range("N136").Value = -250025040

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1056.5
$ws.Range("I5").Value = 562.7778
$ws.Range("J5").Value = 5500
$ws.Range("K5").Value = 1688.3334
$ws.Range("L5").Value = 16500
$ws.Range("M5").Value = -1576.3334
$ws.Range("N5").Value = -16724
$ws.Range("H32").Value = 2108
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 2260
$ws.Range("K32").Value = 4500
$ws.Range("L32").Value = 6780
$ws.Range("M32").Value = -4217
$ws.Range("N32").Value = -7346
$ws.Range("H68").Value = 1629.5056
$ws.Range("I68").Value = 625.55884
$ws.Range("J68").Value = 2250.1272
$ws.Range("K68").Value = 1876.67652
$ws.Range("L68").Value = 6750.3816
$ws.Range("M68").Value = -1065.67652
$ws.Range("N68").Value = -8372.381600000001
$ws.Range("H71").Value = 1629.5056
$ws.Range("I71").Value = 625.55884
$ws.Range("J71").Value = 2250.1272
$ws.Range("K71").Value = 5630.02956
$ws.Range("L71").Value = 20251.1448
$ws.Range("M71").Value = -1574.02956
$ws.Range("N71").Value = -28363.1448
$ws.Range("H107").Value = 826.8387
$ws.Range("J107").Value = 1128.8387
$ws.Range("L107").Value = 3386.5161
$ws.Range("N107").Value = -7226.5161
$ws.Range("H135").Value = 1056.5
$ws.Range("I135").Value = 562.7778
$ws.Range("J135").Value = 5500
$ws.Range("K135").Value = 5065.000199999999
$ws.Range("L135").Value = 49500
$ws.Range("M135").Value = -2530.000199999999
$ws.Range("N135").Value = -54570
$ws.Range("H140").Value = 18519006
$ws.Range("I140").Value = 18519006
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 55557018
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = $null
$ws.Range("N140").Value = -55551838

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4331.24
$ws.Range("I132").Value = 3163.2856
$ws.Range("J132").Value = 5817.727
$ws.Range("K132").Value = 9489.856800000001
$ws.Range("L132").Value = 17453.181
$ws.Range("M132").Value = -6959.856800000001
$ws.Range("N132").Value = -22513.181

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2576.1538
$ws.Range("I40").Value = 1927.7142
$ws.Range("K40").Value = 1927.7142
$ws.Range("M40").Value = -1791.7142
$ws.Range("H61").Value = 333338340
$ws.Range("I61").Value = 500002500
$ws.Range("K61").Value = 500002500
$ws.Range("M61").Value = -500002298
$ws.Range("H68").Value = 1764.6471
$ws.Range("I68").Value = 999.93335
$ws.Range("K68").Value = 999.93335
$ws.Range("M68").Value = -250.93335
$ws.Range("H71").Value = 1764.6471
$ws.Range("I71").Value = 999.93335
$ws.Range("K71").Value = 4999.66675
$ws.Range("M71").Value = -1255.66675
$ws.Range("H113").Value = 333338340
$ws.Range("I113").Value = 500002500
$ws.Range("K113").Value = 500002500
$ws.Range("M113").Value = -500000330

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3450290.8
$ws.Range("I126").Value = 1745.6111
$ws.Range("J126").Value = 9093364
$ws.Range("K126").Value = 5236.8333
$ws.Range("L126").Value = 27280092
$ws.Range("M126").Value = -2766.8333
$ws.Range("N126").Value = -27285032
